# Release-Notes.xlsx update
#
# The folder-inventory scan re-ran. The most recently touched folder,
# "Build A Fabric Real-Time Intelligence Solution in a Day", used to sit at
# the bottom of the "Folder Inventory" table (row 62); it now sorts to the
# top (row 2) with a refreshed "Last Updated" timestamp, and every row that
# used to sit between the header and that entry shifts down by one (rows
# 63 onward are untouched). The "Metadata" and "Summary" sheets get their
# run timestamp / workflow counter / most-recent-update values bumped to
# match.

$wb = $excel.ActiveWorkbook

# ---- Folder Inventory sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Folder Inventory")

# Make room at the top: push rows 2..62 down into 3..63.
$ws.Rows(2).Insert(-4121)  # -4121 = xlShiftDown

# The folder that used to be row 62 is now row 63 -- copy its values up into
# the freshly-opened row 2 (cell by cell; Value2 is the reliable getter here).
for ($col = 1; $col -le 5; $col++) {
    $moved = $ws.Cells.Item(63, $col).Value2
    $ws.Cells.Item(2, $col).Value = $moved
}

# Remove the now-duplicated row that used to hold that data (63), closing
# the gap back up so the table stays at its original 71 data rows.
$ws.Rows(63).Delete(-4162)  # -4162 = xlShiftUp

# Inserting above row 2 pulled in the header row's bold/centered formatting;
# drop it so the relocated row matches the other plain data rows.
$ws.Range("A2:E2").ClearFormats()

# The re-scan refreshed this folder's "Last Updated" timestamp.
$ws.Range("C2").Value = "2025-06-12 15:59:35 +0530"

# ---- Metadata sheet ----------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2025-06-12 10:29:52 UTC"

# "Workflow Run" is stored as text (not a number); force the cell to Text
# format before writing so Excel doesn't silently turn "2" into a numeric 2.
$meta.Range("B5").NumberFormat = "@"
$meta.Range("B5").Value = "2"

# ---- Summary sheet ------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = "2025-06-12 15:59:35 +0530"
